$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("E3").Value = "8FB15"
$ws.Range("F3").Value = "Thông số càng 88734"
$ws.Range("G3").Value = "3454976814"
$ws.Range("I3").Value = 40
$ws.Range("J3").Value = 190
$ws.Range("K3").Value = 60
$ws.Range("L3").Value = 70
$ws.Range("M3").Value = 120

# Row 4
$ws.Range("E4").Value = "FE4P16"
$ws.Range("F4").Value = "Loại động cơ 61557"
$ws.Range("G4").Value = "8671095827"
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 60
$ws.Range("L4").Value = 70
$ws.Range("M4").Value = 130

# Row 5
$ws.Range("E5").Value = "FD30T-16"
$ws.Range("F5").Value = "Số lượng van 17365"
$ws.Range("G5").Value = "7735273602"
$ws.Range("I5").Value = 40
$ws.Range("J5").Value = 170
$ws.Range("K5").Value = 60
$ws.Range("L5").Value = 60
$ws.Range("M5").Value = 110

# Row 6
$ws.Range("E6").Value = "FE4P16"
$ws.Range("F6").Value = "Số lượng van 64553"
$ws.Range("G6").Value = "2340786685"
$ws.Range("I6").Value = 40
$ws.Range("J6").Value = 180
$ws.Range("K6").Value = 80
$ws.Range("L6").Value = 80
$ws.Range("M6").Value = 120

# Row 7
$ws.Range("E7").Value = "FB25-12"
$ws.Range("F7").Value = "Thông số càng 06816"
$ws.Range("G7").Value = "5935655771"
$ws.Range("I7").Value = 40
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 70
$ws.Range("L7").Value = 80
$ws.Range("M7").Value = 120

# Row 8
$ws.Range("E8").Value = "8FD25"
$ws.Range("F8").Value = "Sideshift 13935"
$ws.Range("G8").Value = "6618779189"
$ws.Range("I8").Value = 30
$ws.Range("J8").Value = 180
$ws.Range("K8").Value = 80
$ws.Range("L8").Value = 80
$ws.Range("M8").Value = 130

$ws.Range("F7").Select()
